$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.503.76'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.643.04'
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3819'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3617'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08250'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.234'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.471'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.358'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '1.639.64'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06965'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.606'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").Value = '23.493.56'
$ws.Range("E24").Value = '  +0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.543'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.083'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.274'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").Value = '1.819.92'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("E32").Value = '  +15.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.158'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.564'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.89%  '
$ws.Range("E35").Value = '  +6.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02776'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2517'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08777'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.990'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07039'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7063'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.347'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6545'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.299'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07980'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.85%  '
$ws.Range("E51").Value = '  -0.65%  '
